$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 839.3684
$ws.Range("J32").Value = 803.4706
$ws.Range("L32").Value = 803.4706
$ws.Range("N32").Value = -1455.4706

$ws.Range("H40").Value = 7967.5884
$ws.Range("I40").Value = 6625.3335
$ws.Range("J40").Value = 11189
$ws.Range("K40").Value = 6625.3335
$ws.Range("L40").Value = 11189
$ws.Range("M40").Value = -6450.3335
$ws.Range("N40").Value = -11539

$ws.Range("H53").Value = 25642518
$ws.Range("I53").Value = 66667748
$ws.Range("J53").Value = 1749.25
$ws.Range("K53").Value = 66667748
$ws.Range("L53").Value = 1749.25
$ws.Range("M53").Value = -66667111
$ws.Range("N53").Value = -3023.25

$ws.Range("H64").Value = 6168.8223
$ws.Range("I64").Value = 2800
$ws.Range("J64").Value = 6325.5117
$ws.Range("K64").Value = 2800
$ws.Range("L64").Value = 6325.5117
$ws.Range("N64").Value = -6821.5117
$ws.Range("M64").Value = -2552

$ws.Range("H67").Value = 6168.8223
$ws.Range("I67").Value = 2800
$ws.Range("J67").Value = 6325.5117
$ws.Range("K67").Value = 2800
$ws.Range("L67").Value = 6325.5117
$ws.Range("N67").Value = -8041.5117
$ws.Range("M67").Value = -1942

$ws.Range("H74").Value = 10605
$ws.Range("I74").Value = 8561.23
$ws.Range("K74").Value = 8561.23
$ws.Range("M74").Value = -7625.23

$ws.Range("H77").Value = 10605
$ws.Range("I77").Value = 8561.23
$ws.Range("K77").Value = 42806.14999999999
$ws.Range("M77").Value = -38126.14999999999

$ws.Range("H86").Value = 2927591.5
$ws.Range("J86").Value = 4390148
$ws.Range("L86").Value = 4390148
$ws.Range("N86").Value = -4392394

$ws.Range("H89").Value = 2927591.5
$ws.Range("J89").Value = 4390148
$ws.Range("L89").Value = 21950740
$ws.Range("N89").Value = -21961972

$ws.Range("H112").Value = 3380.9524
$ws.Range("J112").Value = 3380.9524
$ws.Range("L112").Value = 10142.8572
$ws.Range("N112").Value = -12358.8572

$ws.Range("H132").Value = 3619.158
$ws.Range("I132").Value = 3662.5881
$ws.Range("K132").Value = 10987.7643
$ws.Range("M132").Value = -8457.764299999999

$ws.Range("H134").Value = 71017.16
$ws.Range("J134").Value = 71017.16
$ws.Range("L134").Value = 71017.16
$ws.Range("N134").Value = -81157.16

$ws.Range("H138").Value = 6726.5557
$ws.Range("I138").Value = 4091.25
$ws.Range("J138").Value = 7684.8486
$ws.Range("K138").Value = 12273.75
$ws.Range("L138").Value = 23054.5458
$ws.Range("M138").Value = -7133.75
$ws.Range("N138").Value = -33334.5458

$ws.Range("H140").Value = 77499.75
$ws.Range("J140").Value = 77499.75
$ws.Range("L140").Value = 77499.75
$ws.Range("N140").Value = -87859.75

$ws.Range("H141").Value = 4329.0835
$ws.Range("I141").Value = 4329.0835
$ws.Range("K141").Value = 12987.2505
$ws.Range("M141").Value = -7807.250499999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 289.8
$ws.Range("I4").Value = 237.25
$ws.Range("K4").Value = 237.25
$ws.Range("M4").Value = -121.25

$ws.Range("H74").Value = 1504
$ws.Range("J74").Value = 1153.3334
$ws.Range("L74").Value = 1153.3334
$ws.Range("N74").Value = -2901.3334

$ws.Range("H77").Value = 1504
$ws.Range("J77").Value = 1153.3334
$ws.Range("L77").Value = 5766.666999999999
$ws.Range("N77").Value = -14502.667

$ws.Range("H97").Value = 1539.4286
$ws.Range("I97").Value = 1562.75
$ws.Range("K97").Value = 1562.75
$ws.Range("M97").Value = -1066.75

$ws.Range("H102").Value = 2038.3334
$ws.Range("I102").Value = 2270.6
$ws.Range("K102").Value = 2270.6
$ws.Range("M102").Value = -648.5999999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1170.8148
$ws.Range("I94").Value = 1278.762
$ws.Range("J94").Value = 793
$ws.Range("K94").Value = 1278.762
$ws.Range("L94").Value = 793
$ws.Range("M94").Value = -827.7619999999999
$ws.Range("N94").Value = -1695

$ws.Range("H99").Value = 1752.5714
$ws.Range("I99").Value = 1655.7693
$ws.Range("J99").Value = 3011
$ws.Range("K99").Value = 1655.7693
$ws.Range("L99").Value = 3011
$ws.Range("M99").Value = -157.7692999999999
$ws.Range("N99").Value = -6007

$ws.Range("H105").Value = 204201.8
$ws.Range("I105").Value = 204201.8
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 204201.8
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -202454.8
$ws.Range("N105").ClearContents()

$ws.Range("H134").Value = 36539.098
$ws.Range("I134").Value = 4025.4285
$ws.Range("J134").Value = 340000
$ws.Range("K134").Value = 12076.2855
$ws.Range("L134").Value = 1020000
$ws.Range("M134").Value = -9541.2855
$ws.Range("N134").Value = -1025070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 433.9565
$ws.Range("I7").Value = 430.3
$ws.Range("K7").Value = 430.3
$ws.Range("M7").Value = -317.3

$ws.Range("H62").Value = 2115.5
$ws.Range("I62").Value = 1897.5
$ws.Range("J62").Value = 2333.5
$ws.Range("K62").Value = 1897.5
$ws.Range("L62").Value = 2333.5
$ws.Range("M62").Value = -1273.5
$ws.Range("N62").Value = -3581.5

$ws.Range("H65").Value = 2115.5
$ws.Range("I65").Value = 1897.5
$ws.Range("J65").Value = 2333.5
$ws.Range("K65").Value = 9487.5
$ws.Range("L65").Value = 11667.5
$ws.Range("M65").Value = -6367.5
$ws.Range("N65").Value = -17907.5

$ws.Range("H134").Value = 265107.1
$ws.Range("I134").Value = 1998.7428
$ws.Range("J134").Value = 3334704.8
$ws.Range("K134").Value = 5996.2284
$ws.Range("L134").Value = 10004114.4
$ws.Range("M134").Value = -3461.2284
$ws.Range("N134").Value = -10009184.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1064.5
$ws.Range("I23").Value = 435.5
$ws.Range("J23").Value = 1169.3334
$ws.Range("K23").Value = 1306.5
$ws.Range("L23").Value = 3508.0002
$ws.Range("M23").Value = -1071.5
$ws.Range("N23").Value = -3978.0002

$ws.Range("H88").Value = 11999.667
$ws.Range("I88").Value = 13000
$ws.Range("K88").Value = 39000
$ws.Range("M88").Value = -38572

$ws.Range("H91").Value = 11999.667
$ws.Range("I91").Value = 13000
$ws.Range("K91").Value = 39000
$ws.Range("M91").Value = -37518

$ws.Range("H92").Value = 667675.3
$ws.Range("I92").Value = 1250703.9
$ws.Range("K92").Value = 3752111.7
$ws.Range("M92").Value = -3750863.7

$ws.Range("H98").Value = 3988.125
$ws.Range("I98").Value = 5465
$ws.Range("J98").Value = 3495.8333
$ws.Range("K98").Value = 16395
$ws.Range("L98").Value = 10487.4999
$ws.Range("M98").Value = -14897
$ws.Range("N98").Value = -13483.4999

$ws.Range("H131").Value = 2098.4546
$ws.Range("I131").Value = 2120
$ws.Range("J131").Value = 2093.6667
$ws.Range("K131").Value = 6360
$ws.Range("L131").Value = 6281.000100000001
$ws.Range("M131").Value = -1320
$ws.Range("N131").Value = -16361.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 687.1070999999999
$ws.Range("I97").Value = 676.95654
$ws.Range("K97").Value = 676.95654
$ws.Range("M97").Value = -180.95654

$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws.Range("H122").Value = 3799
$ws.Range("I122").Value = 3799
$ws.Range("K122").Value = 11397
$ws.Range("M122").Value = -8947

$ws.Range("H135").Value = 166750000
$ws.Range("J135").Value = 166750000
$ws.Range("L135").Value = 166750000
$ws.Range("N135").Value = -166760140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4290.0303
$ws.Range("I7").Value = 4000.913
$ws.Range("J7").Value = 4955
$ws.Range("K7").Value = 4000.913
$ws.Range("L7").Value = 4955
$ws.Range("M7").Value = -3888.913
$ws.Range("N7").Value = -5179

$ws.Range("H100").Value = 2424.3333
$ws.Range("I100").Value = 2174.25
$ws.Range("J100").Value = 2924.5
$ws.Range("K100").Value = 2174.25
$ws.Range("L100").Value = 2924.5
$ws.Range("M100").Value = -1633.25
$ws.Range("N100").Value = -4006.5

$ws.Range("H126").Value = 4290.0303
$ws.Range("I126").Value = 4000.913
$ws.Range("J126").Value = 4955
$ws.Range("K126").Value = 12002.739
$ws.Range("L126").Value = 14865
$ws.Range("M126").Value = -9532.739
$ws.Range("N126").Value = -19805

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 79642.21000000001
$ws.Range("I62").Value = 257250.25
$ws.Range("K62").Value = 257250.25
$ws.Range("M62").Value = -256626.25

$ws.Range("H65").Value = 79642.21000000001
$ws.Range("I65").Value = 257250.25
$ws.Range("K65").Value = 1286251.25
$ws.Range("M65").Value = -1283131.25

$ws.Range("H113").Value = 280.73334
$ws.Range("I113").Value = 296.91666
$ws.Range("K113").Value = 890.7499799999999
$ws.Range("M113").Value = 1279.25002

$ws.Range("H132").Value = 15868.208
$ws.Range("I132").Value = 2277.6292
$ws.Range("K132").Value = 6832.8876
$ws.Range("M132").Value = -4302.8876

$ws.Range("H135").Value = 105000
$ws.Range("J135").Value = 105000
$ws.Range("L135").Value = 105000
$ws.Range("N135").Value = -115140

$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200

$ws.Range("H139").Value = 54998
$ws.Range("J139").Value = 54998
$ws.Range("L139").Value = 54998
$ws.Range("N139").Value = -65278
